$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 197, pushing the existing
# rows 197:251 down to 200:254 (dimension grows from A1:T251 to A1:T254).
$ws.Rows("197:199").Insert()

# --- New row 197 ---
$ws.Range("A197").Value = 1
$ws.Range("B197").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C197").Value = "Arica y Parinacota"
$ws.Range("D197").Value = 45027
$ws.Range("E197").Value = 15
$ws.Range("F197").Value = "Fruta"
$ws.Range("G197").Value = 100108
$ws.Range("H197").Value = "Tropicales y subtropicales"
$ws.Range("I197").Value = 100108002
$ws.Range("J197").Value = "Mango"
$ws.Range("K197").Value = "Piqueño"
$ws.Range("L197").Value = "Primera"
$ws.Range("M197").Value = 80
$ws.Range("N197").Value = 12000
$ws.Range("O197").Value = 13000
$ws.Range("P197").Value = 12312
$ws.Range("Q197").Value = "$/caja 10 kilos"
$ws.Range("R197").Value = "Región de Arica y Parinacota"
$ws.Range("S197").Value = 1231
$ws.Range("T197").Value = 10

# --- New row 198 ---
$ws.Range("A198").Value = 1
$ws.Range("B198").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C198").Value = "Arica y Parinacota"
$ws.Range("D198").Value = 45027
$ws.Range("E198").Value = 15
$ws.Range("F198").Value = "Fruta"
$ws.Range("G198").Value = 100108
$ws.Range("H198").Value = "Tropicales y subtropicales"
$ws.Range("I198").Value = 100108002
$ws.Range("J198").Value = "Mango"
$ws.Range("K198").Value = "Piqueño"
$ws.Range("L198").Value = "Segunda"
$ws.Range("M198").Value = 180
$ws.Range("N198").Value = 9000
$ws.Range("O198").Value = 10000
$ws.Range("P198").Value = 9500
$ws.Range("Q198").Value = "$/caja 10 kilos"
$ws.Range("R198").Value = "Región de Arica y Parinacota"
$ws.Range("S198").Value = 950
$ws.Range("T198").Value = 10

# --- New row 199 ---
$ws.Range("A199").Value = 1
$ws.Range("B199").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C199").Value = "Arica y Parinacota"
$ws.Range("D199").Value = 45027
$ws.Range("E199").Value = 15
$ws.Range("F199").Value = "Fruta"
$ws.Range("G199").Value = 100108
$ws.Range("H199").Value = "Tropicales y subtropicales"
$ws.Range("I199").Value = 100108002
$ws.Range("J199").Value = "Mango"
$ws.Range("K199").Value = "Piqueño"
$ws.Range("L199").Value = "Tercera"
$ws.Range("M199").Value = 140
$ws.Range("N199").Value = 7000
$ws.Range("O199").Value = 8000
$ws.Range("P199").Value = 7429
$ws.Range("Q199").Value = "$/caja 10 kilos"
$ws.Range("R199").Value = "Región de Arica y Parinacota"
$ws.Range("S199").Value = 743
$ws.Range("T199").Value = 10
